$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.171.22"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "3.800.13"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "598.54"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "170.80"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").Value = "3.796.62"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").Value = "6.53"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "0.0000265"
$ws.Range("E13").Value = "  -4.98%  "
$ws.Range("D14").Value = "36.86"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "4.440.35"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "3.803.62"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "69.170.59"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "18.28"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("D19").Value = "7.10"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "11.11"
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("D22").Value = "471.54"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "0.710"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "84.74"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "0.0000148"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "12.22"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "3.950.95"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").Value = "7.46"
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("D33").Value = "2.25"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").Value = "30.33"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").Value = "9.42"
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("D37").Value = "3.755.58"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").Value = "3.52"
$ws.Range("E39").Value = "  -8.68%  "
$ws.Range("D40").Value = "0.140"
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("E41").Value = "  +1.31%  "
$ws.Range("D42").Value = "5.88"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "0.310"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "1.98"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("D47").Value = "43.79"
$ws.Range("E47").Value = "  +11.58%  "
$ws.Range("D48").Value = "8.66"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").Value = "46.08"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").Value = "403.01"
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").Value = "145.17"
$ws.Range("E51").Value = "  +2.79%  "
